$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1332
$ws1.Range("F3").Value = 2849

# Sheet "全部类型" (all types) - same two events duplicated on rows 3-4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1332
$ws4.Range("F4").Value = 2849
